# Apply the "changes to test and data" commit:
#  - add a new test email value and point the details!C2 cell at it
#  - make "details" the active/selected sheet with selection F6
#  - "indetails" is no longer the selected sheet

$wb = $excel.ActiveWorkbook

$details   = $wb.Worksheets.Item("details")
$indetails = $wb.Worksheets.Item("indetails")

# New test data value (becomes a new shared string, referenced by C2)
$details.Range("C2").Value = "finleaptesto@grr.la"

# "details" becomes the active sheet/tab, with the new selection
$details.Activate()
$details.Range("F6").Select()

# "indetails" (previously the active tab) is no longer tab-selected;
# its own cell selection (D15) is left untouched.
